$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Insert a new row at 152 (shifts existing rows 152:178 down to 153:179),
# copy the formatting from the row above so the new row matches the
# surrounding data rows (style 6 for B, style 1 for C:G).
$ws.Rows("152:152").Insert() | Out-Null
$ws.Range("B151:G151").Copy()
$ws.Range("B152:G152").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New alias entry: UnigetUI / ug / [Application]
# Alias ("ug") is entered before the app name so the shared-string table
# gets "ug" then "UnigetUI", matching how the row was authored.
$ws.Range("D152").Value = "ug"
$ws.Range("C152").Value = "UnigetUI"
$ws.Range("E152").Value = "[Application]"

# Re-enter the SEQUENCE array formula so it spills across the now-larger
# range (COUNTA(C:C)-1 grew by one row) instead of leaving a stale spill.
$ws.Range("B3:B175").FormulaArray = "=SEQUENCE(COUNTA(C:C)-1)"

# Restore the view state (scroll position / selection) to match the saved
# workbook.
$ws.Range("C157").Select() | Out-Null
